$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5", "D6", "D8", "D11", "D12", "D13", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D29", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '74.870.54'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").Value = '2.842.57'
$ws.Range("E3").Value = '  +9.68%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '603.52'
$ws.Range("E5").Value = '  +3.62%  '
$ws.Range("D6").Value = '189.50'
$ws.Range("E6").Value = '  +2.78%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.557'
$ws.Range("E8").Value = '  +4.77%  '
$ws.Range("E9").Value = '  -6.02%  '
$ws.Range("D10").Value = '2.840.77'
$ws.Range("E10").Value = '  +9.70%  '
$ws.Range("D11").Value = '0.162'
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("D12").Value = '0.372'
$ws.Range("E12").Value = '  +3.31%  '
$ws.Range("D13").Value = '4.92'
$ws.Range("E13").Value = '  +2.31%  '
$ws.Range("D14").Value = '3.366.33'
$ws.Range("E14").Value = '  +10.09%  '
$ws.Range("D15").Value = '74.892.28'
$ws.Range("E15").Value = '  +1.12%  '
$ws.Range("D16").Value = '27.56'
$ws.Range("E16").Value = '  +5.17%  '
$ws.Range("D17").Value = '0.0000189'
$ws.Range("E17").Value = '  -1.88%  '
$ws.Range("D18").Value = '2.841.29'
$ws.Range("E18").Value = '  +10.15%  '
$ws.Range("D19").Value = '9.18'
$ws.Range("E19").Value = '  +8.09%  '
$ws.Range("D20").Value = '12.48'
$ws.Range("E20").Value = '  +6.11%  '
$ws.Range("D21").Value = '377.85'
$ws.Range("E21").Value = '  +3.39%  '
$ws.Range("D22").Value = '2.29'
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("D23").Value = '4.14'
$ws.Range("E23").Value = '  +1.58%  '
$ws.Range("E24").Value = '  -0.21%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = '71.12'
$ws.Range("E26").Value = '  +2.89%  '
$ws.Range("D27").Value = '4.25'
$ws.Range("E27").Value = '  +1.45%  '
$ws.Range("D28").Value = '2.995.77'
$ws.Range("E28").Value = '  +10.45%  '
$ws.Range("D29").Value = '9.67'
$ws.Range("E29").Value = '  +4.96%  '
$ws.Range("E30").Value = '  +11.83%  '
$ws.Range("E31").Value = '  +0.12%  '
$ws.Range("D32").Value = '532.31'
$ws.Range("E32").Value = '  +6.73%  '
$ws.Range("D33").Value = '1.41'
$ws.Range("E33").Value = '  +5.76%  '
$ws.Range("D34").Value = '7.97'
$ws.Range("E34").Value = '  +0.79%  '
$ws.Range("D35").Value = '1.82'
$ws.Range("E35").Value = '  +6.64%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").Value = '0.121'
$ws.Range("E37").Value = '  +2.59%  '
$ws.Range("D38").Value = '20.22'
$ws.Range("E38").Value = '  +5.33%  '
$ws.Range("D39").Value = '162.33'
$ws.Range("E39").Value = '  +1.40%  '
$ws.Range("D40").Value = '19.29'
$ws.Range("E40").Value = '  -0.52%  '
$ws.Range("D41").Value = '183.93'
$ws.Range("E41").Value = '  +23.50%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").Value = '5.14'
$ws.Range("E43").Value = '  +5.75%  '
$ws.Range("D44").Value = '0.342'
$ws.Range("E44").Value = '  +7.39%  '
$ws.Range("D45").Value = '1.69'
$ws.Range("E45").Value = '  +1.39%  '
$ws.Range("D46").Value = '1.26'
$ws.Range("E46").Value = '  +8.63%  '
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").Value = '2.39'
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").Value = '39.76'
$ws.Range("E48").Value = '  +1.65%  '
$ws.Range("D49").Value = '0.0853'
$ws.Range("E49").Value = '  +5.03%  '
$ws.Range("D50").Value = '0.575'
$ws.Range("E50").Value = '  +10.46%  '
$ws.Range("D51").Value = '3.76'
$ws.Range("E51").Value = '  +4.39%  '
